# correção do parâmetro spb para ano base 2021
$wb = $excel.ActiveWorkbook

$wsMeses = $wb.Worksheets.Item("meses")
$wsMeses.Range("B5").Value = 0.3

$wsFator = $wb.Worksheets.Item("fator")
$wsFator.Select()
$wsFator.Range("B2").Select()

$wsMeses.Select()
$wsMeses.Range("B6").Select()

$wb.Save()
